# Applies the "Quantum Biology" -> "Chemistry" rewrite described by the
# commit diff, using the Word COM object model (Find/Replace + Range.Text).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Unveiling the Microscopic Realm: A Journey into Quantum Biology",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unveiling the Enigma of Chemistry: A Journey into the Realm of Matter and Transformations",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Author name
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Emily Taylor",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dr. Avery Donovan",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3) Author e-mail
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "etaylor@biodiscovery.org",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "chemistrydr.donovan@highschool.academy",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) Main body paragraph (the long paragraph with the manual line breaks)
# ---------------------------------------------------------------------
$body = $d.Paragraphs.Item(5).Range
$body.MoveEnd(1, -1) | Out-Null
$body.Text = (
    "In the vast tapestry of sciences, chemistry stands as a beacon of discovery, " +
    "illuminating the composition of matter and the intricate dance of transformations it undergoes." +
    " As we embark on this captivating journey into the realm of chemistry, we will unravel the secrets " +
    "of the atom, witness the symphony of chemical reactions, and explore the profound implications of " +
    "chemistry in medicine, industry, and everyday life.`v`v" +
    "Chemistry unveils the enigmatic world of particles, revealing the intricate dance of atoms, ions, and molecules." +
    " We will delve into the depths of the periodic table, deciphering the patterns and properties that govern the elements." +
    " From the fiery brilliance of lithium to the noble elegance of helium, each element holds a unique story, waiting to be explored.`v`v" +
    "The macroscopic world is a stage on which chemistry plays a transformative role." +
    " Chemical reactions, like choreographed ballets, orchestrate spectacular displays of color, energy, and matter." +
    " From the explosive combustion of fuels to the gentle rusting of iron, chemical reactions shape our world in countless ways.`v`v" +
    "Beyond the theoretical realm, chemistry finds boundless applications in medicine, industry, and everyday life." +
    " From life-saving drugs to durable materials, chemistry touches every aspect of our existence." +
    " We will uncover the myriad ways in which chemistry contributes to our health, well-being, and technological advancements."
)

# ---------------------------------------------------------------------
# 5) "Summary" heading paragraph is unchanged.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 6) Summary body paragraph
# ---------------------------------------------------------------------
$summary = $d.Paragraphs.Item(7).Range
$summary.MoveEnd(1, -1) | Out-Null
$summary.Text = (
    "In this essay, we embarked on a captivating expedition into the world of chemistry, delving into the " +
    "microscopic realm of particles and exploring the enchanting choreography of chemical reactions." +
    " We discovered the profound applications of chemistry in medicine, industry, and everyday life." +
    " Chemistry, with its ability to unravel the secrets of matter and orchestrate transformations, stands as a " +
    "pillar of scientific knowledge, enriching our understanding of the universe and empowering us to create a better future."
)

# ---------------------------------------------------------------------
# 7) Trailing empty paragraph added at the end of the document.
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0) | Out-Null
$endRange.InsertParagraphAfter() | Out-Null
